{"js": "// Load the paragraphs so we can find the second paragraph\n// (\"Word doc for anything that happens up ...\").\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// 1) Merge the three runs (split apart by the now-removed grammar-check\n//    <w:proofErr> markers) of the second paragraph back into one run/one\n//    sentence with normal (non-preserved) spacing.\nconst introPara = paragraphs.items[1];\nintroPara.clear();\nintroPara.insertText(\n  \"Word doc for anything that happens up until the player makes their choice for who to go on a date with\",\n  Word.InsertLocation.start\n);\n\n// 2) Append the new plan paragraphs after the intro paragraph, in order.\nconst newParagraphText = [\n  \"\",\n  \"Im not certain yet but I have a vague idea of how it goes\",\n  \"First, game starts with you choosing your name (already in file)\",\n  \"Then, cuts to y/n monologuing to themselves about how they gotta do errands (they hate errands).\",\n  \"Starts by going to the grocery store.\",\n  \"After some tough grocerying, y/n leans up against a shelf to relax\",\n  \"Shelf fall\",\n  \"Goo crushed\",\n  \"Apologizes made, gets paper towels to clean him up\",\n  \"Buys paper towels\",\n  \"Heads to doctors afterwards for appointment\",\n  \"Waiting in hallway, bothered by weird clown girl\",\n  \"Kinda chill convo, may or may not be charming\",\n  \"Dove floats out of their sleeve off to nowhere\",\n  \"Gotta chase it down together\",\n  \"Shes apologetic, offers to take u out for coffee\",\n  \"Head to post office to finish up\",\n  \"Guess whos there\",\n  \"Goo.\",\n  \"Have a chat, apologize, laugh about it, and ngl this goo guy is a lil cute\",\n  \"Get his #\",\n  \"Go home\",\n  \"Before sleeping, chose who you wanna date\",\n  \"Something something \\u201c___ is so pathetic, barely able to function, massive pushover, and living a insane life. Yet, they are so lovely, I think I can fix them.\\u201d\",\n  \"Cut to each route.\",\n];\n\nlet anchor = introPara;\nfor (const text of newParagraphText) {\n  anchor = anchor.insertParagraph(text, Word.InsertLocation.after);\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) The second paragraph (\"Word doc for anything that ...\") used to be split\n#    across three runs around a grammar-check <w:proofErr> pair\n#    (\"happens up\"). Replace the whole paragraph's text (but not its\n#    paragraph mark) in one shot so it collapses back down to a single run\n#    / a single sentence, with no leftover proofErr markers.\n$introPara = $d.Paragraphs(2)\n$introRange = $d.Range($introPara.Range.Start, $introPara.Range.End)\n$introRange.Text = \"Word doc for anything that happens up until the player makes their choice for who to go on a date with\"\n\n# 2) Append the new planning paragraphs right after the intro paragraph.\n$leftDQuote = [string][char]0x201C\n$rightDQuote = [string][char]0x201D\n\n$newLines = @(\n  \"Im not certain yet but I have a vague idea of how it goes\",\n  \"First, game starts with you choosing your name (already in file)\",\n  \"Then, cuts to y/n monologuing to themselves about how they gotta do errands (they hate errands).\",\n  \"Starts by going to the grocery store.\",\n  \"After some tough grocerying, y/n leans up against a shelf to relax\",\n  \"Shelf fall\",\n  \"Goo crushed\",\n  \"Apologizes made, gets paper towels to clean him up\",\n  \"Buys paper towels\",\n  \"Heads to doctors afterwards for appointment\",\n  \"Waiting in hallway, bothered by weird clown girl\",\n  \"Kinda chill convo, may or may not be charming\",\n  \"Dove floats out of their sleeve off to nowhere\",\n  \"Gotta chase it down together\",\n  \"Shes apologetic, offers to take u out for coffee\",\n  \"Head to post office to finish up\",\n  \"Guess whos there\",\n  \"Goo.\",\n  \"Have a chat, apologize, laugh about it, and ngl this goo guy is a lil cute\",\n  \"Get his #\",\n  \"Go home\",\n  \"Before sleeping, chose who you wanna date\",\n  (\"Something something \" + $leftDQuote + \"___ is so pathetic, barely able to function, massive pushover, and living a insane life. Yet, they are so lovely, I think I can fix them.\" + $rightDQuote),\n  \"Cut to each route.\"\n)\n\n$introPara = $d.Paragraphs(2)\n$insertPoint = $d.Range($introPara.Range.End, $introPara.Range.End)\n\n# Blank separator paragraph, inserted the same way Word itself creates an\n# empty paragraph (no text run content at all).\n$insertPoint.InsertParagraphAfter()\n\n# Re-fetch the (now third) paragraph -- the freshly inserted blank one --\n# and drop the rest of the plan text after it in a single bulk insert. The\n# leading \"`r\" starts a brand new paragraph instead of typing into the\n# blank one.\n$blankPara = $d.Paragraphs(3)\n$afterBlank = $d.Range($blankPara.Range.End, $blankPara.Range.End)\n$afterBlank.InsertAfter(\"`r\" + ($newLines -join \"`r\"))\n"}
